$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 14:22"

# Apply updated COVID-19 country statistics and re-sorted rank swaps
$ws.Range("B4").Value = 5842062
$ws.Range("C4").Value = 634
$ws.Range("E4").Value = 2513791
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 180191
$ws.Range("B6").Value = 3050326
$ws.Range("C6").Value = 6890
$ws.Range("D6").Value = 2282358
$ws.Range("E6").Value = 711085
$ws.Range("G6").Value = 37
$ws.Range("H6").Value = 56883
$ws.Range("A18").Value = "Banglades"
$ws.Range("B18").Value = 294598
$ws.Range("C18").Value = 1973
$ws.Range("D18").Value = 179091
$ws.Range("E18").Value = 111566
$ws.Range("G18").Value = 34
$ws.Range("H18").Value = 3941
$ws.Range("A19").Value = "Pakistan"
$ws.Range("B19").Value = 292765
$ws.Range("C19").Value = 591
$ws.Range("D19").Value = 275836
$ws.Range("E19").Value = 10694
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 6235
$ws.Range("B23").Value = 233912
$ws.Range("C23").Value = 55
$ws.Range("E23").Value = 15631
$ws.Range("B28").Value = 117008
$ws.Range("C28").Value = 243
$ws.Range("D28").Value = 113808
$ws.Range("E28").Value = 3007
$ws.Range("B41").Value = 80528
$ws.Range("C41").Value = 571
$ws.Range("D41").Value = 72307
$ws.Range("E41").Value = 7706
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 515
$ws.Range("B42").Value = 78505
$ws.Range("C42").Value = 961
$ws.Range("D42").Value = 35287
$ws.Range("E42").Value = 39946
$ws.Range("G42").Value = 39
$ws.Range("H42").Value = 3272
$ws.Range("B59").Value = 39903
$ws.Range("C59").Value = 276
$ws.Range("E59").Value = 3802
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 2001
$ws.Range("B68").Value = 31935
$ws.Range("C68").Value = 818
$ws.Range("D68").Value = 18631
$ws.Range("E68").Value = 13155
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 149
$ws.Range("B80").Value = 16317
$ws.Range("C80").Value = 78
$ws.Range("D80").Value = 14180
$ws.Range("E80").Value = 1515
$ws.Range("B82").Value = 14327
$ws.Range("C82").Value = 50
$ws.Range("D82").Value = 13355
$ws.Range("E82").Value = 794
$ws.Range("B85").Value = 12949
$ws.Range("C85").Value = 99
$ws.Range("D85").Value = 8455
$ws.Range("E85").Value = 4225
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 269
$ws.Range("A99").Value = "Croacia"
$ws.Range("B99").Value = 8175
$ws.Range("C99").Value = 275
$ws.Range("D99").Value = 5801
$ws.Range("E99").Value = 2203
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 171
$ws.Range("A100").Value = "Haiti"
$ws.Range("B100").Value = 8050
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 5447
$ws.Range("E100").Value = 2407
$ws.Range("H100").Value = 196
$ws.Range("A101").Value = "Finlandia"
$ws.Range("B101").Value = 7920
$ws.Range("C101").Value = 14
$ws.Range("D101").Value = 7100
$ws.Range("E101").Value = 486
$ws.Range("H101").Value = 334
$ws.Range("B133").Value = 2272
$ws.Range("C133").Value = 7
$ws.Range("E133").Value = 185
$ws.Range("B139").Value = 2064
$ws.Range("C139").Value = 6
$ws.Range("D139").Value = 1939
$ws.Range("E139").Value = 115
$ws.Range("A159").Value = "Vietnam"
$ws.Range("B159").Value = 1016
$ws.Range("C159").Value = 2
$ws.Range("D159").Value = 563
$ws.Range("E159").Value = 426
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 27
$ws.Range("A160").Value = "Lesoto"
$ws.Range("B160").Value = 1015
$ws.Range("D160").Value = 472
$ws.Range("E160").Value = 513
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 30
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
